$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row-level D/E updates (rows 2-47) ---
$ws.Range("D2").Value = '26.772.36'
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").Value = '1.648.74'
$ws.Range("E3").Value = '  +1.41%  '
$ws.Range("E4").Value = '  +0.57%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.39'
$ws.Range("E5").Value = '  +1.70%  '
$ws.Range("E6").Value = '  +1.38%  '
$ws.Range("E7").Value = '  +0.42%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.251'
$ws.Range("E8").Value = '  +1.87%  '
$ws.Range("E9").Value = '  +0.80%  '
$ws.Range("E10").Value = '  +2.61%  '
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("D12").Value = '1.878.76'
$ws.Range("E12").Value = '  +1.44%  '
$ws.Range("D13").Value = '1.656.12'
$ws.Range("E13").Value = '  +2.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.19'
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("E15").Value = '  +2.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.45'
$ws.Range("E16").Value = '  +0.86%  '
$ws.Range("D17").Value = '26.776.21'
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").Value = '0.0₃0744'
$ws.Range("E18").Value = '  +0.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '218.68'
$ws.Range("E19").Value = '  +2.34%  '
$ws.Range("E20").Value = '  +0.41%  '
$ws.Range("E21").Value = '  +2.07%  '
$ws.Range("E22").Value = '  +0.68%  '
$ws.Range("E23").Value = '  +15.77%  '
$ws.Range("E24").Value = '  +2.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.73'
$ws.Range("E25").Value = '  -1.25%  '
$ws.Range("E26").Value = '  +0.53%  '
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("E28").Value = '  +4.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.75'
$ws.Range("E29").Value = '  +1.71%  '
$ws.Range("E30").Value = '  +1.88%  '
$ws.Range("E31").Value = '  +2.04%  '
$ws.Range("E32").Value = '  +1.16%  '
$ws.Range("E33").Value = '  +2.42%  '
$ws.Range("D34").Value = '1.281.56'
$ws.Range("E34").Value = '  +4.58%  '
$ws.Range("E35").Value = '  +3.66%  '
$ws.Range("E36").Value = '  +2.35%  '
$ws.Range("E37").Value = '  +3.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.538'
$ws.Range("E38").Value = '  +6.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.829'
$ws.Range("E39").Value = '  +4.44%  '
$ws.Range("E40").Value = '  +0.49%  '
$ws.Range("E41").Value = '  +2.70%  '
$ws.Range("E42").Value = '  -1.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.47'
$ws.Range("E43").Value = '  +2.60%  '
$ws.Range("D44").Value = '1.789.19'
$ws.Range("E44").Value = '  +1.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.08'
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.77'
$ws.Range("E46").Value = '  +9.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.60'
$ws.Range("E47").Value = '  +1.96%  '

# --- Rows 48-51: BabyDogeCoin inserted, list shifts down, Mantle drops off ---
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0103'
$ws.Range("E48").Value = '  +2.51%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0516'
$ws.Range("E49").Value = '  +1.17%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.76'
$ws.Range("E50").Value = '  +3.93%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0970'
$ws.Range("E51").Value = '  +2.20%  '